$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 values (order matters for shared-string indices: A8, B8, F8, D8)
$ws.Range("A8").Value = "MEA (Midwest Economic Association)"
$ws.Range("B8").Value = "St Louis"

# C8: date, reuse the same number-format style as the other date cells (e.g. C4)
$ws.Range("C4").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = "10/8/2018"

# F8: comment text, reuse the same style as A8 (left aligned, bordered)
$ws.Range("A8").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = "Abstract submission for members only. Member fees: $20/year, abstract submission fee: $20, conference registration fee: $45, for non-members $135."

# D8: date-range text with its own larger, distinctly colored font
$ws.Range("D8").Value = "March 15-17, 2019"
$ws.Range("D8").Style = "Normal"
$ws.Range("D8").Font.Size = 16
$ws.Range("D8").Font.Name = "Trebuchet MS"
$ws.Range("D8").Font.Color = 3026478

# Update the sheet view to match the saved selection
$ws.Range("D4").Select() | Out-Null
